# Update cryptos list figures (price + 1h volume change) per scheduled data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '40.087.21'
$ws.Range("E2").Value = '  +2.83%  '
$ws.Range("D3").Value = '2.235.84'
$ws.Range("E3").Value = '  +0.97%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''294.95'
$ws.Range("E5").Value = '  -0.27%  '
$ws.Range("E6").Value = '  +8.08%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").Value = '''0.473'
$ws.Range("E9").Value = '  +3.75%  '
$ws.Range("D10").Value = '''31.15'
$ws.Range("E10").Value = '  +11.54%  '
$ws.Range("D11").Value = '''0.0793'
$ws.Range("E11").Value = '  +2.88%  '
$ws.Range("D12").Value = '''47.03'
$ws.Range("E12").Value = '  +2.53%  '
$ws.Range("D13").Value = '''0.109'
$ws.Range("E13").Value = '  +1.14%  '
$ws.Range("D14").Value = '''6.47'
$ws.Range("E14").Value = '  +6.39%  '
$ws.Range("D15").Value = '2.579.77'
$ws.Range("E15").Value = '  +0.85%  '
$ws.Range("D16").Value = '''14.18'
$ws.Range("E16").Value = '  +2.32%  '
$ws.Range("D17").Value = '2.180.76'
$ws.Range("E17").Value = '  -2.21%  '
$ws.Range("D18").Value = '''0.730'
$ws.Range("E18").Value = '  +3.16%  '
$ws.Range("D19").Value = '39.993.92'
$ws.Range("E19").Value = '  +2.83%  '
$ws.Range("D20").Value = '0.0₃0890'
$ws.Range("E20").Value = '  +3.84%  '
$ws.Range("D21").Value = '''5.81'
$ws.Range("E21").Value = '  +2.23%  '
$ws.Range("D22").Value = '''10.92'
$ws.Range("E22").Value = '  +12.11%  '
$ws.Range("D23").Value = '''65.50'
$ws.Range("E23").Value = '  +1.58%  '
$ws.Range("D24").Value = '''235.33'
$ws.Range("E24").Value = '  +4.86%  '
$ws.Range("E25").Value = '  +0.06%  '
$ws.Range("D27").Value = '''1.85'
$ws.Range("E27").Value = '  +6.29%  '
$ws.Range("D28").Value = '''22.85'
$ws.Range("E28").Value = '  +2.82%  '
$ws.Range("E29").Value = '  +3.35%  '
$ws.Range("D30").Value = '''9.24'
$ws.Range("E30").Value = '  +4.14%  '
$ws.Range("D31").Value = '''33.44'
$ws.Range("E31").Value = '  +7.90%  '
$ws.Range("D32").Value = '''152.34'
$ws.Range("E32").Value = '  +2.61%  '
$ws.Range("D34").Value = '''4.90'
$ws.Range("E34").Value = '  +3.72%  '
$ws.Range("D35").Value = '''0.0720'
$ws.Range("E35").Value = '  +5.16%  '
$ws.Range("E36").Value = '  +2.29%  '
$ws.Range("D37").Value = '''16.44'
$ws.Range("E37").Value = '  +15.52%  '
$ws.Range("D38").Value = '''0.112'
$ws.Range("E38").Value = '  +3.36%  '
$ws.Range("E39").Value = '  +5.68%  '
$ws.Range("D40").Value = '''2.72'
$ws.Range("E40").Value = '  +3.70%  '
$ws.Range("D41").Value = '''1.70'
$ws.Range("E41").Value = '  +7.53%  '
$ws.Range("E42").Value = '  +6.85%  '
$ws.Range("D43").Value = '2.045.63'
$ws.Range("E43").Value = '  +7.83%  '
$ws.Range("D44").Value = '''2.23'
$ws.Range("E44").Value = '  +7.39%  '
$ws.Range("E45").Value = '  +6.84%  '
$ws.Range("D46").Value = '''10.00'
$ws.Range("E46").Value = '  +13.61%  '
$ws.Range("E47").Value = '  +1.78%  '
$ws.Range("E48").Value = '  +2.79%  '
$ws.Range("D49").Value = '2.451.12'
$ws.Range("E49").Value = '  +1.03%  '
$ws.Range("D50").Value = '''71.08'
$ws.Range("E50").Value = '  +1.94%  '
$ws.Range("D51").Value = '''1.46'
$ws.Range("E51").Value = '  +15.13%  '
